$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 729.8
$ws.Range("I55").Value = 75
$ws.Range("K55").Value = 75
$ws.Range("M55").Value = 139
$ws.Range("H107").Value = 668.7778
$ws.Range("I107").Value = 526.25
$ws.Range("J107").Value = 709.5
$ws.Range("K107").Value = 526.25
$ws.Range("L107").Value = 709.5
$ws.Range("M107").Value = 1393.75
$ws.Range("N107").Value = -4549.5
$ws.Range("H113").Value = 2968.0557
$ws.Range("I113").Value = 2439.4443
$ws.Range("K113").Value = 2439.4443
$ws.Range("M113").Value = 814.5556999999999
$ws.Range("H125").Value = 143728.42
$ws.Range("I125").Value = 200799.8
$ws.Range("J125").Value = 1050
$ws.Range("K125").Value = 1807198.2
$ws.Range("L125").Value = 9450
$ws.Range("M125").Value = -1804738.2
$ws.Range("N125").Value = -14370
$ws.Range("H132").Value = 1803.55
$ws.Range("I132").Value = 1698.4324
$ws.Range("K132").Value = 5095.2972
$ws.Range("M132").Value = -2565.2972
$ws.Range("H138").Value = 2262.225
$ws.Range("I138").Value = 1137.4546
$ws.Range("J138").Value = 3636.9443
$ws.Range("K138").Value = 3412.3638
$ws.Range("L138").Value = 10910.8329
$ws.Range("M138").Value = 1727.6362
$ws.Range("N138").Value = -21190.8329

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 23706.334
$ws.Range("J37").Value = 24447.6
$ws.Range("L37").Value = 24447.6
$ws.Range("N37").Value = -24993.6
$ws.Range("H44").Value = 26983.334
$ws.Range("I44").Value = 1000
$ws.Range("J44").Value = 39975
$ws.Range("K44").Value = 1000
$ws.Range("L44").Value = 39975
$ws.Range("M44").Value = -512
$ws.Range("N44").Value = -40951
$ws.Range("H55").Value = 25825
$ws.Range("I55").Value = 6700
$ws.Range("J55").Value = 44950
$ws.Range("K55").Value = 6700
$ws.Range("L55").Value = 44950
$ws.Range("M55").Value = -6385
$ws.Range("N55").Value = -45580
$ws.Range("H61").Value = 2797.6
$ws.Range("I61").Value = 3710.2856
$ws.Range("J61").Value = 2306.1538
$ws.Range("K61").Value = 3710.2856
$ws.Range("L61").Value = 2306.1538
$ws.Range("M61").Value = -3498.2856
$ws.Range("N61").Value = -2730.1538
$ws.Range("H74").Value = 1113.1333
$ws.Range("I74").Value = 1145
$ws.Range("J74").Value = 1085.25
$ws.Range("K74").Value = 1145
$ws.Range("L74").Value = 1085.25
$ws.Range("M74").Value = -271
$ws.Range("N74").Value = -2833.25
$ws.Range("H77").Value = 1113.1333
$ws.Range("I77").Value = 1145
$ws.Range("J77").Value = 1085.25
$ws.Range("K77").Value = 5725
$ws.Range("L77").Value = 5426.25
$ws.Range("M77").Value = -1357
$ws.Range("N77").Value = -14162.25
$ws.Range("H123").Value = 24272
$ws.Range("J123").Value = 24272
$ws.Range("L123").Value = 24272
$ws.Range("N123").Value = -34072
$ws.Range("H136").Value = 2797.6
$ws.Range("I136").Value = 3710.2856
$ws.Range("J136").Value = 2306.1538
$ws.Range("K136").Value = 11130.8568
$ws.Range("L136").Value = 6918.4614
$ws.Range("M136").Value = -8580.856800000001
$ws.Range("N136").Value = -12018.4614

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H62").Value = 29824.967
$ws.Range("J62").Value = 30152.467
$ws.Range("L62").Value = 30152.467
$ws.Range("N62").Value = -31524.467
$ws.Range("H65").Value = 29824.967
$ws.Range("J65").Value = 30152.467
$ws.Range("L65").Value = 90457.401
$ws.Range("N65").Value = -97321.401
$ws.Range("H80").Value = 2469300.2
$ws.Range("I80").Value = 18518568
$ws.Range("J80").Value = 182.15384
$ws.Range("K80").Value = 18518568
$ws.Range("L80").Value = 182.15384
$ws.Range("M80").Value = -18517570
$ws.Range("N80").Value = -2178.15384
$ws.Range("H83").Value = 2469300.2
$ws.Range("I83").Value = 18518568
$ws.Range("J83").Value = 182.15384
$ws.Range("K83").Value = 92592840
$ws.Range("L83").Value = 910.7692
$ws.Range("M83").Value = -92587848
$ws.Range("N83").Value = -10894.7692
$ws.Range("H134").Value = 3129.125
$ws.Range("I134").Value = 2773.5334
$ws.Range("J134").Value = 3721.7778
$ws.Range("K134").Value = 8320.600199999999
$ws.Range("L134").Value = 11165.3334
$ws.Range("M134").Value = -5785.600199999999
$ws.Range("N134").Value = -16235.3334

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9998.666999999999
$ws.Range("J4").Value = 9998.666999999999
$ws.Range("L4").Value = 9998.666999999999
$ws.Range("N4").Value = -10222.667
$ws.Range("H31").Value = 3194.4285
$ws.Range("I31").Value = 1611.091
$ws.Range("K31").Value = 1611.091
$ws.Range("M31").Value = -1316.091
$ws.Range("H34").Value = 3194.4285
$ws.Range("I34").Value = 1611.091
$ws.Range("K34").Value = 1611.091
$ws.Range("M34").Value = -1409.091
$ws.Range("H58").Value = 1612682.6
$ws.Range("I58").Value = 2180782
$ws.Range("J58").Value = 3067.3333
$ws.Range("K58").Value = 2180782
$ws.Range("L58").Value = 3067.3333
$ws.Range("M58").Value = -2180579
$ws.Range("N58").Value = -3473.3333
$ws.Range("H132").Value = 277202.8
$ws.Range("I132").Value = 398470.72
$ws.Range("J132").Value = 2328.9333
$ws.Range("K132").Value = 1195412.16
$ws.Range("L132").Value = 6986.7999
$ws.Range("M132").Value = -1192882.16
$ws.Range("N132").Value = -12046.7999
$ws.Range("H136").Value = 1612682.6
$ws.Range("I136").Value = 2180782
$ws.Range("J136").Value = 3067.3333
$ws.Range("K136").Value = 6542346
$ws.Range("L136").Value = 9201.999899999999
$ws.Range("M136").Value = -6539796
$ws.Range("N136").Value = -14301.9999

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2331.3333
$ws.Range("I4").Value = 797.6
$ws.Range("J4").Value = 10000
$ws.Range("K4").Value = 2392.8
$ws.Range("L4").Value = 30000
$ws.Range("M4").Value = -2280.8
$ws.Range("N4").Value = -30224
$ws.Range("H25").Value = 1885.4445
$ws.Range("I25").Value = 533.3333
$ws.Range("J25").Value = 4589.6665
$ws.Range("K25").Value = 1599.9999
$ws.Range("L25").Value = 13768.9995
$ws.Range("M25").Value = -1430.9999
$ws.Range("N25").Value = -14106.9995
$ws.Range("H30").Value = 1885.4445
$ws.Range("I30").Value = 533.3333
$ws.Range("J30").Value = 4589.6665
$ws.Range("K30").Value = 1599.9999
$ws.Range("L30").Value = 13768.9995
$ws.Range("M30").Value = -1497.9999
$ws.Range("N30").Value = -13972.9995
$ws.Range("H32").Value = 1837.7778
$ws.Range("J32").Value = 1837.7778
$ws.Range("L32").Value = 5513.3334
$ws.Range("N32").Value = -6079.3334
$ws.Range("H38").Value = 129.11111
$ws.Range("I38").Value = 120
$ws.Range("J38").Value = 147.33333
$ws.Range("K38").Value = 360
$ws.Range("L38").Value = 441.99999
$ws.Range("M38").Value = -13
$ws.Range("N38").Value = -1135.99999
$ws.Range("H117").Value = 63678.5
$ws.Range("I117").Value = 514.5
$ws.Range("J117").Value = 84733.164
$ws.Range("K117").Value = 1543.5
$ws.Range("L117").Value = 254199.492
$ws.Range("M117").Value = 1898.5
$ws.Range("N117").Value = -261083.492
$ws.Range("H131").Value = 13700159
$ws.Range("I131").Value = 394
$ws.Range("J131").Value = 14707495
$ws.Range("K131").Value = 1182
$ws.Range("L131").Value = 44122485
$ws.Range("M131").Value = 3858
$ws.Range("N131").Value = -44132565
$ws.Range("H132").Value = 1931.84
$ws.Range("I132").Value = 1535.7646
$ws.Range("J132").Value = 2773.5
$ws.Range("K132").Value = 13821.8814
$ws.Range("L132").Value = 24961.5
$ws.Range("M132").Value = -11291.8814
$ws.Range("N132").Value = -30021.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2043.2424
$ws.Range("I102").Value = 1713.0385
$ws.Range("J102").Value = 3269.7144
$ws.Range("K102").Value = 1713.0385
$ws.Range("L102").Value = 3269.7144
$ws.Range("M102").Value = -91.03850000000011
$ws.Range("N102").Value = -6513.7144
$ws.Range("H109").Value = 9275.6
$ws.Range("J109").Value = 9275.6
$ws.Range("L109").Value = 9275.6
$ws.Range("N109").Value = -11355.6
$ws.Range("H113").Value = 1910.2354
$ws.Range("J113").Value = 2436.75
$ws.Range("L113").Value = 2436.75
$ws.Range("N113").Value = -6776.75

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1330.875
$ws.Range("J22").Value = 999
$ws.Range("L22").Value = 999
$ws.Range("N22").Value = -1589
$ws.Range("H27").Value = 1330.875
$ws.Range("J27").Value = 999
$ws.Range("L27").Value = 999
$ws.Range("N27").Value = -1213
$ws.Range("H122").Value = 16670629
$ws.Range("I122").Value = 19233872
$ws.Range("J122").Value = 14290475
$ws.Range("K122").Value = 57701616
$ws.Range("L122").Value = 42871425
$ws.Range("M122").Value = -57699166
$ws.Range("N122").Value = -42876325
$ws.Range("H136").Value = 24635922
$ws.Range("I136").Value = 34484050
$ws.Range("J136").Value = 836284.5600000001
$ws.Range("K136").Value = 103452150
$ws.Range("L136").Value = 2508853.68
$ws.Range("M136").Value = -103449600
$ws.Range("N136").Value = -2513953.68

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 357334
$ws.Range("J5").Value = 36001
$ws.Range("L5").Value = 36001
$ws.Range("N5").Value = -36225
$ws.Range("H123").Value = 23787.387
$ws.Range("J123").Value = 23787.387
$ws.Range("L123").Value = 23787.387
$ws.Range("N123").Value = -33587.387
